# Actualizar WEEK, TASK, ultimas semanas
# Actualizacion de scripts WEEK, TASK del equipo de las ultimas semanas
#
# Updates the hours-worked figures for the last couple of logged weeks on
# "Hoja1" (K30 / K31). The dependent running-total column (L) is a chain of
# formulas (L_n = K_n + L_{n-1}), so editing just these two source cells
# ripples the new totals through the rest of the table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("K30").Value = 4
$ws.Range("K31").Value = 8

# Leave the view roughly where the author last left it: scrolled down to the
# most recently edited rows, with the just-touched row (32) selected.
$ws.Range("A19").Select()
$ws.Range("A32:XFD32").Select()
